$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.570.06'
$ws.Range("E2").Value = '  +6.94%  '

$ws.Range("D3").Value = '1.733.60'
$ws.Range("E3").Value = '  +4.54%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9951'
$ws.Range("E4").Value = '  -0.57%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '333.94'
$ws.Range("E5").Value = '  +5.52%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9954'
$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3716'
$ws.Range("E7").Value = '  +2.35%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.50'
$ws.Range("E8").Value = '  +5.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3386'
$ws.Range("E9").Value = '  +3.66%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.206'
$ws.Range("E10").Value = '  +5.92%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07545'
$ws.Range("E11").Value = '  +6.85%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9929'
$ws.Range("E12").Value = '  -0.48%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.421'
$ws.Range("E13").Value = '  +6.25%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.52'
$ws.Range("E14").Value = '  +4.88%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.009'
$ws.Range("E15").Value = '  +5.92%  '

$ws.Range("D16").Value = '1.725.73'
$ws.Range("E16").Value = '  +3.66%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001095'
$ws.Range("E17").Value = '  +4.39%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06696'
$ws.Range("E18").Value = '  +1.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '83.43'
$ws.Range("E19").Value = '  +5.32%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9948'
$ws.Range("E20").Value = '  -0.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.79'
$ws.Range("E21").Value = '  +6.62%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.183'
$ws.Range("E22").Value = '  +4.55%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.19'
$ws.Range("E23").Value = '  +4.92%  '

$ws.Range("D24").Value = '26.488.85'
$ws.Range("E24").Value = '  +6.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.464'
$ws.Range("E25").Value = '  +1.18%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.513'
$ws.Range("E26").Value = '  +4.61%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.437'
$ws.Range("E27").Value = '  +17.24%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '152.25'
$ws.Range("E28").Value = '  +2.33%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.60'
$ws.Range("E29").Value = '  +5.22%  '

$ws.Range("D30").Value = '1.913.04'
$ws.Range("E30").Value = '  +3.49%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '131.59'
$ws.Range("E31").Value = '  +4.52%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.111'
$ws.Range("E32").Value = '  +0.86%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.131'
$ws.Range("E33").Value = '  +5.26%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08603'
$ws.Range("E34").Value = '  +1.98%  '

$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.19'
$ws.Range("E35").Value = '  +7.46%  '

$ws.Range("B36").Value = 'WEMIXTOKEN'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.710'
$ws.Range("E36").Value = '  +1.94%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.471'
$ws.Range("E37").Value = '  +5.06%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02358'
$ws.Range("E38").Value = '  +5.58%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06363'
$ws.Range("E39").Value = '  +5.37%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.722'
$ws.Range("E40").Value = '  +6.09%  '

$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2176'
$ws.Range("E41").Value = '  +5.10%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.239'
$ws.Range("E42").Value = '  -4.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6265'
$ws.Range("E43").Value = '  +5.68%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.41'
$ws.Range("E44").Value = '  +13.19%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9947'
$ws.Range("E45").Value = '  -0.14%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.900'
$ws.Range("E46").Value = '  +2.41%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6062'
$ws.Range("E47").Value = '  +7.41%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.62'
$ws.Range("E48").Value = '  +3.59%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.069'
$ws.Range("E49").Value = '  +5.98%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07337'
$ws.Range("E50").Value = '  +4.56%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '77.87'
$ws.Range("E51").Value = '  +3.97%  '
